# Update existing row 2 and append two new rows of crawled data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new crawl result for permissionspolicy.xyz
$ws.Range("A2").Value = "https://permissionspolicy.xyz/"
$ws.Range("B2").Value = $true
$ws.Range("C2").Value = $true
$ws.Range("D2").Value = $true
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = "['geolocation']"
$ws.Range("G2").Value = "['permissionspolicy.cdn-gamma.com']"

# Add row 3 for google.com
$ws.Range("A3").Value = "https://google.com/"
$ws.Range("B3").Value = $false
$ws.Range("C3").Value = $false
$ws.Range("D3").Value = $false
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = "[]"
$ws.Range("G3").Value = "[]"

# Add row 4 for facebook.com
$ws.Range("A4").Value = "https://facebook.com/"
$ws.Range("B4").Value = $false
$ws.Range("C4").Value = $false
$ws.Range("D4").Value = $false
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = "[]"
$ws.Range("G4").Value = "[]"
